$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1113.4054
$ws.Range("I28").Value = 352.57693
$ws.Range("J28").Value = 2911.7273
$ws.Range("K28").Value = 352.57693
$ws.Range("L28").Value = 2911.7273
$ws.Range("M28").Value = 132.42307
$ws.Range("N28").Value = -3881.7273
$ws.Range("H40").Value = 8085.5713
$ws.Range("J40").Value = 4728.769
$ws.Range("L40").Value = 4728.769
$ws.Range("N40").Value = -5078.769
$ws.Range("H86").Value = 5492.909
$ws.Range("I86").Value = 6450.3076
$ws.Range("K86").Value = 6450.3076
$ws.Range("M86").Value = -5327.3076
$ws.Range("H89").Value = 5492.909
$ws.Range("I89").Value = 6450.3076
$ws.Range("K89").Value = 32251.538
$ws.Range("M89").Value = -26635.538
$ws.Range("H96").Value = 1133.1666
$ws.Range("I96").Value = 1399.5
$ws.Range("J96").Value = 1000
$ws.Range("K96").Value = 4198.5
$ws.Range("L96").Value = 3000
$ws.Range("M96").Value = -2825.5
$ws.Range("N96").Value = -5746
$ws.Range("H111").Value = 6539108
$ws.Range("J111").Value = 2848
$ws.Range("L111").Value = 8544
$ws.Range("N111").Value = -14678
$ws.Range("H113").Value = 5500
$ws.Range("I113").Value = 4500
$ws.Range("J113").Value = 5600
$ws.Range("K113").Value = 4500
$ws.Range("L113").Value = 5600
$ws.Range("M113").Value = -1246
$ws.Range("N113").Value = -12108
$ws.Range("H116").Value = 11438.8
$ws.Range("I116").Value = 4836.2856
$ws.Range("J116").Value = 17216
$ws.Range("K116").Value = 4836.2856
$ws.Range("L116").Value = 17216
$ws.Range("M116").Value = -1394.2856
$ws.Range("N116").Value = -24100
$ws.Range("H118").Value = 66667164
$ws.Range("I118").Value = 76923490
$ws.Range("J118").Value = 1049.5
$ws.Range("K118").Value = 230770470
$ws.Range("L118").Value = 3148.5
$ws.Range("M118").Value = -230768813
$ws.Range("N118").Value = -6462.5
$ws.Range("H138").Value = 2881.026
$ws.Range("I138").Value = 1719.1786
$ws.Range("J138").Value = 3544.9387
$ws.Range("K138").Value = 5157.5358
$ws.Range("L138").Value = 10634.8161
$ws.Range("M138").Value = -17.53579999999965
$ws.Range("N138").Value = -20914.8161

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2428.95
$ws.Range("I32").Value = 2428.95
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 2428.95
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -2141.95
$ws.Range("N32").ClearContents()
$ws.Range("H61").Value = 5605.854
$ws.Range("I61").Value = 6667
$ws.Range("J61").Value = 3837.2778
$ws.Range("K61").Value = 6667
$ws.Range("L61").Value = 3837.2778
$ws.Range("M61").Value = -6455
$ws.Range("N61").Value = -4261.2778
$ws.Range("H96").Value = 36793.57
$ws.Range("J96").Value = 36793.57
$ws.Range("L96").Value = 36793.57
$ws.Range("N96").Value = -42285.57
$ws.Range("H97").Value = 1598581.6
$ws.Range("I97").Value = 2097707.2
$ws.Range("J97").Value = 1379.8
$ws.Range("K97").Value = 2097707.2
$ws.Range("L97").Value = 1379.8
$ws.Range("M97").Value = -2097211.2
$ws.Range("N97").Value = -2371.8
$ws.Range("H110").Value = 1463399.6
$ws.Range("I110").Value = 1544421.8
$ws.Range("K110").Value = 1544421.8
$ws.Range("M110").Value = -1542376.8
$ws.Range("H136").Value = 5605.854
$ws.Range("I136").Value = 6667
$ws.Range("J136").Value = 3837.2778
$ws.Range("K136").Value = 20001
$ws.Range("L136").Value = 11511.8334
$ws.Range("M136").Value = -17451
$ws.Range("N136").Value = -16611.8334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()
$ws.Range("H80").Value = 457.41666
$ws.Range("I80").Value = 538.6
$ws.Range("J80").Value = 399.42856
$ws.Range("K80").Value = 538.6
$ws.Range("L80").Value = 399.42856
$ws.Range("M80").Value = 459.4
$ws.Range("N80").Value = -2395.42856
$ws.Range("H83").Value = 457.41666
$ws.Range("I83").Value = 538.6
$ws.Range("J83").Value = 399.42856
$ws.Range("K83").Value = 2693
$ws.Range("L83").Value = 1997.1428
$ws.Range("M83").Value = 2299
$ws.Range("N83").Value = -11981.1428
$ws.Range("H94").Value = 2224995.5
$ws.Range("I94").Value = 2632617
$ws.Range("J94").Value = 12193.286
$ws.Range("K94").Value = 2632617
$ws.Range("L94").Value = 12193.286
$ws.Range("M94").Value = -2632166
$ws.Range("N94").Value = -13095.286
$ws.Range("H99").Value = 14288035
$ws.Range("I99").Value = 15875150
$ws.Range("J99").Value = 3997
$ws.Range("K99").Value = 15875150
$ws.Range("L99").Value = 3997
$ws.Range("M99").Value = -15873652
$ws.Range("N99").Value = -6993
$ws.Range("H105").Value = 3678884.2
$ws.Range("I105").Value = 3908627
$ws.Range("K105").Value = 3908627
$ws.Range("M105").Value = -3906880
$ws.Range("H107").Value = 7938293.5
$ws.Range("I107").Value = 11905788
$ws.Range("K107").Value = 11905788
$ws.Range("M107").Value = -11903868

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 42298.25
$ws.Range("J51").Value = 42298.25
$ws.Range("L51").Value = 42298.25
$ws.Range("N51").Value = -43770.25
$ws.Range("H59").Value = 53750
$ws.Range("J59").Value = 57500
$ws.Range("L59").Value = 57500
$ws.Range("N59").Value = -59790
$ws.Range("H61").Value = 42298.25
$ws.Range("J61").Value = 42298.25
$ws.Range("L61").Value = 42298.25
$ws.Range("N61").Value = -42994.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H61").Value = 135.8
$ws.Range("I61").Value = 84.5
$ws.Range("J61").Value = 170
$ws.Range("K61").Value = 253.5
$ws.Range("L61").Value = 510
$ws.Range("M61").Value = -38.5
$ws.Range("N61").Value = -940
$ws.Range("H107").Value = 583.6667
$ws.Range("I107").Value = 762.4
$ws.Range("K107").Value = 2287.2
$ws.Range("M107").Value = -367.1999999999998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 3687.6453
$ws.Range("I2").Value = 484.16
$ws.Range("J2").Value = 17035.5
$ws.Range("K2").Value = 484.16
$ws.Range("L2").Value = 17035.5
$ws.Range("M2").Value = -371.16
$ws.Range("N2").Value = -17261.5
$ws.Range("H12").Value = 49996.332
$ws.Range("J12").Value = 49996.332
$ws.Range("L12").Value = 49996.332
$ws.Range("N12").Value = -50276.332
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").ClearContents()
$ws.Range("H97").Value = 680874.9399999999
$ws.Range("I97").Value = 882312.3
$ws.Range("K97").Value = 882312.3
$ws.Range("M97").Value = -881816.3
$ws.Range("H107").Value = 8438.691999999999
$ws.Range("I107").Value = 17259
$ws.Range("J107").Value = 878.4286
$ws.Range("K107").Value = 17259
$ws.Range("L107").Value = 878.4286
$ws.Range("M107").Value = -15339
$ws.Range("N107").Value = -4718.4286

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4490.5
$ws.Range("I7").Value = 3369.923
$ws.Range("J7").Value = 5814.8184
$ws.Range("K7").Value = 3369.923
$ws.Range("L7").Value = 5814.8184
$ws.Range("M7").Value = -3257.923
$ws.Range("N7").Value = -6038.8184
$ws.Range("H26").Value = 5333.3335
$ws.Range("I26").Value = 2000
$ws.Range("K26").Value = 2000
$ws.Range("M26").Value = -1705
$ws.Range("H122").Value = 6309.5557
$ws.Range("I122").Value = 4578.1
$ws.Range("K122").Value = 13734.3
$ws.Range("M122").Value = -11284.3
$ws.Range("H126").Value = 4490.5
$ws.Range("I126").Value = 3369.923
$ws.Range("J126").Value = 5814.8184
$ws.Range("K126").Value = 10109.769
$ws.Range("L126").Value = 17444.4552
$ws.Range("M126").Value = -7639.769
$ws.Range("N126").Value = -22384.4552
$ws.Range("H136").Value = 84874.72
$ws.Range("J136").Value = 7666.5835
$ws.Range("L136").Value = 22999.7505
$ws.Range("N136").Value = -28099.7505

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 20999.666
$ws.Range("J31").Value = 20999.666
$ws.Range("L31").Value = 20999.666
$ws.Range("N31").Value = -21695.666
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()
$ws.Range("H62").Value = 12743.208
$ws.Range("J62").Value = 9131.727999999999
$ws.Range("L62").Value = 9131.727999999999
$ws.Range("N62").Value = -10379.728
$ws.Range("H65").Value = 12743.208
$ws.Range("J65").Value = 9131.727999999999
$ws.Range("L65").Value = 45658.64
$ws.Range("N65").Value = -51898.64
$ws.Range("H136").Value = 3931.9666
$ws.Range("I136").Value = 4341.8477
$ws.Range("J136").Value = 2585.2144
$ws.Range("K136").Value = 13025.5431
$ws.Range("L136").Value = 7755.6432
$ws.Range("M136").Value = -10475.5431
$ws.Range("N136").Value = -12855.6432
